$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Change the "Run Mode" column (C) for rows 11-31 from "Yes" to "No"
$ws.Range("C11:C31").Value = "No"

# 2) Grow the AutoFilter range from A1:F28 to A1:F31.
#    (Done before row 32 is populated below - otherwise the engine
#    auto-expands the filter to cover the newly-added row too.)
$ws.AutoFilterMode = $false
$ws.Range("A1:F31").AutoFilter()

# 3) Update the hidden _FilterDatabase defined name to match the new AutoFilter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Automation Tests!_FilterDatabase") {
        $n.RefersTo = "='Automation Tests'!`$A`$1:`$F`$31"
    }
}

# 4) Append the new "Field_Mapping_TC001" test case as row 32
# (cells are written in this order so new shared strings land at the same
#  indices - 86,87,88 - as in the target workbook)
$ws.Cells.Item(32, 4).Value = "Load Validated Successfully"
$ws.Cells.Item(32, 1).Value = "Field_Mapping_TC001"
$ws.Cells.Item(32, 2).Value = "1) Enter valid user id and Password and click Login button in Scoular Shipper User.`n2) Click on Add New Load button.`n3) Enter valid details in required field.`n4) Select any value from drop down fields in Field Mapping at right hand side.`n5) Click on Save button.`n6) Select load and click Edit button.`n7) Enter valid details in required field.`n8) Select any value from drop down fields in Field Mapping at right hand side.`n9) Click on Save button.`n10)Select any Load and click Delete button."
$ws.Cells.Item(32, 2).WrapText = $true
$ws.Cells.Item(32, 3).Value = "Yes"
$ws.Rows.Item(32).RowHeight = 195

# 5) Move the selection / scrolled view to the newly added row
$ws.Range("B31").Select()
